$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A9").Value = 18
$ws.Range("B9").Value = "Alex F"
$ws.Range("C9").Value = "alex@exmaple.come"
$ws.Range("D9").Value = "active"
$ws.Range("E9").Value = "none"
$ws.Range("F9").Font.Name = "Calibri"
$ws.Range("G9").Font.Name = "Calibri"
$ws.Range("H9").Value = "2024-12-09 03:16:11"
$ws.Range("I9").Value = "2024-12-09 03:16:11"
$ws.Range("J9").Value = 0
$ws.Range("K9").Formula = '=TEXT(0,"0.0%")'
$ws.Range("K9").Copy()
$ws.Range("K9").PasteSpecial(-4163)
